$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new YouTube link value in D3
$ws.Range("D3").Value = "https://youtu.be/zxTC0YBY2RY"

# Widen column D to fit the new content (target stored width 64.21875 chars;
# the COM width setter here quantizes to whole-pixel steps, so 63.3 is the
# input that lands closest to the target after round-tripping)
$ws.Columns.Item(4).ColumnWidth = 63.3

# Update the active selection to D3
$ws.Range("D3").Select()
